$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two input values driving the wheel analysis recalculation.
# A4: Wheel radius (in) - confirmed from 3 to 4
$ws.Range("A4").Value = 4
# A10: Support count - increased from 6 to 12 (added wheels to BOM)
$ws.Range("A10").Value = 12

# Recalculate the workbook so all dependent formulas refresh their cached values.
$excel.CalculateFullRebuild()

# Update the current selection to reflect where the user left off reviewing the sheet.
$ws.Activate()
$ws.Range("J13").Select()

$wb.Save()
